# TC06_Bento_MultiFilter_Arm-Diagnosis-TumorSize-PRStatus-EndocrineTher.xlsx
# "updated bento tc as per bento perf data availability"
#
# The Neo4j/Web Cypher queries stored on the "startup" sheet filter cases by
# tp.endocrine_therapy_type. Flip the filter value from "OFS" to "Tam" in
# every query cell that references it (columns B and C of rows 2-4), then
# leave the sheet scrolled/selected the way the author left it (cell D2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldFilter = 'endocrine_therapy_type IN ["OFS"]'
$newFilter = 'endocrine_therapy_type IN ["Tam"]'

for ($r = 2; $r -le 4; $r++) {
    for ($c = 2; $c -le 3; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val.Contains($oldFilter)) {
            $cell.Value = $val.Replace($oldFilter, $newFilter)
        }
    }
}

# Restore the view: scrolled so row 2 is at the top, with D2 selected
# (matches topLeftCell="A2" / selection activeCell="D2" sqref="D2").
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("D2").Select()
